$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: rain_effect for stage "1b" changes from "Y" to text "2" (quote-prefixed text)
$ws.Range("D3").Value = "'2"

# N3 / N5: new measure-anchor references ("allow none" / "allow weight")
$ws.Range("N3").Value = "#measure-000000080029218"
$ws.Range("N5").Value = "#measure-000000087253635,#measure-000000129816215"

# Update the active selection to N5, matching the saved view state
$ws.Range("N5").Select()
